# chore: update Sheets via scheduled runner
# Refreshes the cached market-price / leve-profit figures (columns H-N) on the
# per-job "Ixion_Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with
# newly pulled values. Rows whose leve reward has no HQ data end up with the
# HQ profit cell (N) cleared entirely instead of holding a stale number.

$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3500
$ws.Range("I74").Value = 4071.4285
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 4071.4285
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -3135.4285
$ws.Range("N74").Value = -4872

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 11598.857
$ws.Range("I76").Value = 19747.666
$ws.Range("J76").Value = 5487.25
$ws.Range("K76").Value = 19747.666
$ws.Range("L76").Value = 5487.25
$ws.Range("M76").Value = -19432.666
$ws.Range("N76").Value = -6117.25

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3500
$ws.Range("I77").Value = 4071.4285
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 20357.1425
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -15677.1425
$ws.Range("N77").Value = -24360

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 11598.857
$ws.Range("I79").Value = 19747.666
$ws.Range("J79").Value = 5487.25
$ws.Range("K79").Value = 19747.666
$ws.Range("L79").Value = 5487.25
$ws.Range("M79").Value = -18655.666
$ws.Range("N79").Value = -7671.25

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3335068.2
$ws.Range("I113").Value = 5130051.5
$ws.Range("K113").Value = 5130051.5
$ws.Range("M113").Value = -5126797.5

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 35351.668
$ws.Range("I116").Value = 51252.5
$ws.Range("J116").Value = 3550
$ws.Range("K116").Value = 51252.5
$ws.Range("L116").Value = 3550
$ws.Range("M116").Value = -47810.5
$ws.Range("N116").Value = -10434

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6991.616
$ws.Range("I32").Value = 4841.7534
$ws.Range("J32").Value = 25384.889
$ws.Range("K32").Value = 4841.7534
$ws.Range("L32").Value = 25384.889
$ws.Range("M32").Value = -4554.7534
$ws.Range("N32").Value = -25958.889

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 125004860
$ws.Range("I63").Value = 142861280
$ws.Range("J63").Value = 9900
$ws.Range("K63").Value = 142861280
$ws.Range("L63").Value = 9900
$ws.Range("M63").Value = -142860594
$ws.Range("N63").Value = -11272

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 125004860
$ws.Range("I66").Value = 142861280
$ws.Range("J66").Value = 9900
$ws.Range("K66").Value = 714306400
$ws.Range("L66").Value = 49500
$ws.Range("M66").Value = -714302968
$ws.Range("N66").Value = -56364

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 193.33333
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 180
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 180
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -526

# BSM row 32
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 17278.357
$ws.Range("I105").Value = 21071.545
$ws.Range("J105").Value = 3370
$ws.Range("K105").Value = 21071.545
$ws.Range("L105").Value = 3370
$ws.Range("M105").Value = -19324.545
$ws.Range("N105").Value = -6864

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4429.2856
$ws.Range("I62").Value = 5706.6665
$ws.Range("J62").Value = 3471.25
$ws.Range("K62").Value = 5706.6665
$ws.Range("L62").Value = 3471.25
$ws.Range("M62").Value = -5082.6665
$ws.Range("N62").Value = -4719.25

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4429.2856
$ws.Range("I65").Value = 5706.6665
$ws.Range("J65").Value = 3471.25
$ws.Range("K65").Value = 28533.3325
$ws.Range("L65").Value = 17356.25
$ws.Range("M65").Value = -25413.3325
$ws.Range("N65").Value = -23596.25

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7819910
$ws.Range("I99").Value = 15243.143
$ws.Range("J99").Value = 13890207
$ws.Range("K99").Value = 15243.143
$ws.Range("L99").Value = 13890207
$ws.Range("M99").Value = -13745.143
$ws.Range("N99").Value = -13893203

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7819910
$ws.Range("I126").Value = 15243.143
$ws.Range("J126").Value = 13890207
$ws.Range("K126").Value = 45729.429
$ws.Range("L126").Value = 41670621
$ws.Range("M126").Value = -43259.429
$ws.Range("N126").Value = -41675561

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 7143079
$ws.Range("I97").Value = 7692493
$ws.Range("K97").Value = 23077479
$ws.Range("M97").Value = -23076983

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1409418.5
$ws.Range("I131").Value = 14286168
$ws.Range("J131").Value = 1023.9844
$ws.Range("K131").Value = 42858504
$ws.Range("L131").Value = 3071.9532
$ws.Range("M131").Value = -42853464
$ws.Range("N131").Value = -13151.9532

# GSM row 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6450
$ws.Range("I70").Value = 6500.1816
$ws.Range("J70").Value = 6266
$ws.Range("K70").Value = 6500.1816
$ws.Range("L70").Value = 6266
$ws.Range("M70").Value = -6230.1816
$ws.Range("N70").Value = -6806

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6450
$ws.Range("I73").Value = 6500.1816
$ws.Range("J73").Value = 6266
$ws.Range("K73").Value = 6500.1816
$ws.Range("L73").Value = 6266
$ws.Range("M73").Value = -5564.1816
$ws.Range("N73").Value = -8138

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12246.1
$ws.Range("I80").Value = 18625.834
$ws.Range("J80").Value = 2676.5
$ws.Range("K80").Value = 18625.834
$ws.Range("L80").Value = 2676.5
$ws.Range("M80").Value = -17627.834
$ws.Range("N80").Value = -4672.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 12246.1
$ws.Range("I83").Value = 18625.834
$ws.Range("J83").Value = 2676.5
$ws.Range("K83").Value = 93129.17
$ws.Range("L83").Value = 13382.5
$ws.Range("M83").Value = -88137.17
$ws.Range("N83").Value = -23366.5

# LTW row 47
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# LTW row 52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1162.3726
$ws.Range("I136").Value = 683.29034
$ws.Range("J136").Value = 1904.95
$ws.Range("K136").Value = 2049.87102
$ws.Range("L136").Value = 5714.85
$ws.Range("M136").Value = 500.12898
$ws.Range("N136").Value = -10814.85
